$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '21.695.50'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '1.538.46'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').Value = "'289.42"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.12%  '
$ws.Range('D7').Value = "'0.3906"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +3.98%  '
$ws.Range('D8').Value = "'0.3164"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.11%  '
$ws.Range('D9').Value = "'42.84"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.14%  '
$ws.Range('D10').Value = "'0.07164"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('D11').Value = "'1.052"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -6.00%  '
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = "'5.608"
$ws.Range('D13').ClearFormats()
$ws.Range('D14').Value = "'18.53"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.71%  '
$ws.Range('D15').Value = "'6.606"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.73%  '
$ws.Range('D16').Value = '1.546.02'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('D17').Value = "'0.00001095"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.59%  '
$ws.Range('D18').Value = "'0.06569"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range('D19').Value = "'82.92"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.34%  '
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('E21').Value = '  -4.28%  '
$ws.Range('D22').Value = "'15.31"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.92%  '
$ws.Range('D23').Value = "'10.84"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -5.60%  '
$ws.Range('D24').Value = "'2.378"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +5.43%  '
$ws.Range('D25').Value = '21.715.15'
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('D26').Value = "'2.348"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -6.16%  '
$ws.Range('D27').Value = "'147.51"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('D28').Value = "'18.35"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.39%  '
$ws.Range('D29').Value = "'4.838"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '1.718.46'
$ws.Range('E30').Value = '  -0.52%  '
$ws.Range('D31').Value = "'117.08"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.36%  '
$ws.Range('D32').Value = "'0.9590"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -14.48%  '
$ws.Range('D33').Value = "'5.853"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.91%  '
$ws.Range('D34').Value = "'0.08167"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('D35').Value = "'8.706"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -5.76%  '
$ws.Range('D36').Value = "'0.06052"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.54%  '
$ws.Range('D37').Value = "'5.090"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.61%  '
$ws.Range('D38').Value = "'0.02193"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.73%  '
$ws.Range('D39').Value = "'0.2028"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.72%  '
$ws.Range('D40').Value = "'1.175"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.99%  '
$ws.Range('D41').Value = "'1.420"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -13.68%  '
$ws.Range('D43').Value = "'10.59"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.34%  '
$ws.Range('D44').Value = "'0.5707"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.48%  '
$ws.Range('D45').Value = "'3.734"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('D46').Value = "'12.95"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.57%  '
$ws.Range('D47').Value = "'0.5454"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.46%  '
$ws.Range('D48').Value = "'1.159"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('D49').Value = "'115.72"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.52%  '
$ws.Range('D50').Value = "'1.858"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.36%  '
$ws.Range('D51').Value = "'0.06696"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.86%  '
